$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Cypher query stored in the "FilesTab" row (B4): the File Type
# column and the Breed column are removed from the RETURN clause.
$newQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n MATCH (samp:sample)-->(c) `n WHERE samp.specific_sample_pathology IN [`"Melanoma`"]  `nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newQuery

# Update the saved selection so it points at the cell that was edited.
$ws.Range("B4").Select()
